$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lm = "91697550"
$title = "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Series A1 TCL"
$price = "3.699.00"

for ($r = 67; $r -le 73; $r++) {
    # Leading apostrophe forces text storage for the purely-numeric LM code,
    # matching the original column A cells (all stored as text in the sheet).
    $ws.Cells.Item($r, 1).Value = "'" + $lm
    $ws.Cells.Item($r, 2).Value = $title
    $ws.Cells.Item($r, 3).Value = $price
}
